# "vu đã tahy đổi" - add a new member (Trinh Dinh Vu) to NHOM 6, right
# after Do Thi Hong Vy (row 32), and scroll/select down to the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 33: name / phone / email (email styled + hyperlinked like F32)
$ws.Range("D33").Value = "Trịnh Đình Vũ"
$ws.Range("E33").Value = 961973654
$ws.Range("F33").Value = "trinhvu21899@gmail.com"

$ws.Hyperlinks.Add($ws.Range("F33"), "mailto:trinhvu21899@gmail.com") | Out-Null

# Match the look of the existing hyperlink cell (F32) exactly.
$ws.Range("F33").Style = $ws.Range("F32").Style

# Move the view/selection down onto the freshly added row, like the author did.
$ws.Range("F33").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 30
